{"js": "// Applies the \"Companions Description\" edits:\n// 1. Remus section: \"find a women\" -> \"find a woman\"\n// 2. Remus section: \"has do endure\" -> \"has to endure\"\n// 3. Cecilia section: \"She don't likes to talk too much ... this young women wants\"\n//                   -> \"She doesn't like to talk much ... this young woman wants\"\n// 4. Avius section: append two brand-new paragraphs (\"Turn:\" and the follow-up\n//    paragraph about absorbing other Gifted) after the final paragraph of the\n//    document.\n\nconst body = context.document.body;\n\n// --- Edit 1: \"find a women and start a family\" -> \"find a woman and start a family\"\nconst search1 = body.search(\"find a women and start a family\", { matchCase: true });\nsearch1.load(\"text\");\nawait context.sync();\nif (search1.items.length > 0) {\n  search1.items[0].insertText(\"find a woman and start a family\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2: \"ensure nobody has do endure such a cruel fate.\" -> \"... has to endure ...\"\nconst search2 = body.search(\"ensure nobody has do endure such a cruel fate.\", { matchCase: true });\nsearch2.load(\"text\");\nawait context.sync();\nif (search2.items.length > 0) {\n  search2.items[0].insertText(\"ensure nobody has to endure such a cruel fate.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 3: fix \"don't likes ... too much\" / \"this young women\" in Cecilia's description\nconst search3 = body.search(\n  \"She don\\u2019t likes to talk too much and does not trust anyone but herself. Whatever this young women wants, she\",\n  { matchCase: true }\n);\nsearch3.load(\"text\");\nawait context.sync();\nif (search3.items.length > 0) {\n  search3.items[0].insertText(\n    \"She doesn\\u2019t like to talk much and does not trust anyone but herself. Whatever this young woman wants, she\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Edit 4: append the two new \"Turn\" paragraphs at the very end of the document\nconst lastParagraph = body.paragraphs.getLast();\nlastParagraph.load(\"text\");\nawait context.sync();\n\n// First new paragraph: bold \"Turn\" label followed by regular body text.\nconst turnParagraph = lastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nconst turnLabel = turnParagraph.insertText(\"Turn\", Word.InsertLocation.end);\nturnLabel.font.bold = true;\nawait context.sync();\n\nconst turnBody = turnParagraph.insertText(\n  \": Avius is an extremely powerful Gifted. However, that was not always the case. He found out during his battles with the Ridonyan Empire, that a Giften with enough knowledge and willpower, can absorb the power of another Gifted. This is done by consuming their blood and imbuing oneself with the power of the other Gifted. This results in one of the Gifted, namely the one who\\u2019s power is going to be absorbed, dying.\",\n  Word.InsertLocation.end\n);\nturnBody.font.bold = false;\nawait context.sync();\n\n// Second new paragraph: continuation, entirely regular (non-bold) text.\nconst absorbParagraph = turnParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nconst absorbBody = absorbParagraph.insertText(\n  \"Throughout his battles, Avius has been consuming the power and lifeforce of countless Blinded and even some Gifted he saw as lacking. His ultimate goal is the absorption of every Gifted in the world. He believes that when he achieves his goal, he turns into a god himself and that he can \\u201crewrite\\u201d the world to fit his image. With no pain and injustice.\",\n  Word.InsertLocation.end\n);\nabsorbBody.font.bold = false;\nawait context.sync();\n", "ps1": "# Applies the \"Companions Description\" edits:\n# 1. Remus section: \"find a women\" -> \"find a woman\"\n# 2. Remus section: \"has do endure\" -> \"has to endure\"\n# 3. Cecilia section: \"She don't likes to talk too much ... this young women wants\"\n#                   -> \"She doesn't like to talk much ... this young woman wants\"\n# 4. Avius section: append two brand-new paragraphs (\"Turn:\" and the follow-up\n#    paragraph about absorbing other Gifted) after the final paragraph of the\n#    document.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"find a women and start a family\" -> \"find a woman and start a family\"\n$find1 = $d.Content.Find\n$find1.Execute(\n    \"find a women and start a family\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"find a woman and start a family\",\n    1\n) | Out-Null\n\n# --- Edit 2: \"ensure nobody has do endure such a cruel fate.\" -> \"... has to endure ...\"\n$find2 = $d.Content.Find\n$find2.Execute(\n    \"ensure nobody has do endure such a cruel fate.\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"ensure nobody has to endure such a cruel fate.\",\n    1\n) | Out-Null\n\n# --- Edit 3: fix \"don't likes ... too much\" / \"this young women\" in Cecilia's description\n$apostrophe = [char]0x2019\n$findText3 = \"She don\" + $apostrophe + \"t likes to talk too much and does not trust anyone but herself. Whatever this young women wants, she\"\n$replaceText3 = \"She doesn\" + $apostrophe + \"t like to talk much and does not trust anyone but herself. Whatever this young woman wants, she\"\n$find3 = $d.Content.Find\n$find3.Execute(\n    $findText3,\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    $replaceText3,\n    1\n) | Out-Null\n\n# --- Edit 4: append the two new \"Turn\" paragraphs at the very end of the document\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n# First new paragraph: bold \"Turn\" label followed by regular body text.\n$turnParagraph = $d.Paragraphs.Last\n$turnStart = $turnParagraph.Range.Start\n\n$labelRange = $d.Range($turnStart, $turnStart)\n$labelRange.InsertAfter(\"Turn\")\n$labelRange.Font.Bold = 1\n\n$bodyStart = $turnStart + 4\n$rightsQuoteBuilder = [char]0x2019\n$turnBodyText = \": Avius is an extremely powerful Gifted. However, that was not always the case. He found out during his battles with the Ridonyan Empire, that a Giften with enough knowledge and willpower, can absorb the power of another Gifted. This is done by consuming their blood and imbuing oneself with the power of the other Gifted. This results in one of the Gifted, namely the one who\" + $rightsQuoteBuilder + \"s power is going to be absorbed, dying.\"\n$bodyRange = $d.Range($bodyStart, $bodyStart)\n$bodyRange.InsertAfter($turnBodyText)\n$bodyRange.Font.Bold = 0\n\n# Second new paragraph: continuation, entirely regular (non-bold) text.\n$turnParagraph = $d.Paragraphs.Last\n$turnParagraph.Range.InsertParagraphAfter()\n$absorbParagraph = $d.Paragraphs.Last\n$leftQuote = [char]0x201C\n$rightQuote = [char]0x201D\n$absorbText = \"Throughout his battles, Avius has been consuming the power and lifeforce of countless Blinded and even some Gifted he saw as lacking. His ultimate goal is the absorption of every Gifted in the world. He believes that when he achieves his goal, he turns into a god himself and that he can \" + $leftQuote + \"rewrite\" + $rightQuote + \" the world to fit his image. With no pain and injustice.\"\n$absorbRange = $d.Range($absorbParagraph.Range.Start, $absorbParagraph.Range.Start)\n$absorbRange.InsertAfter($absorbText)\n$absorbRange.Font.Bold = 0\n"}
